# Apply the Jan 6 2024 cryptos-list refresh (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin names, links, percentages, and prices that are not
# ambiguous with a plain number) - assigned directly as strings.
$textUpdates = [ordered]@{
    'D2' = '44.061.28'
    'E2' = '  +0.04%  '
    'D3' = '2.239.89'
    'E3' = '  -0.17%  '
    'E4' = '  +0.09%  '
    'E5' = '  -3.45%  '
    'E6' = '  -5.56%  '
    'E7' = '  -0.44%  '
    'E8' = '  +0.12%  '
    'E9' = '  -3.12%  '
    'E10' = '  -4.98%  '
    'E11' = '  -2.10%  '
    'E12' = '  -3.12%  '
    'E13' = '  -0.59%  '
    'D14' = '2.581.33'
    'E14' = '  -0.20%  '
    'D15' = '2.241.99'
    'E15' = '  -0.28%  '
    'E16' = '  -2.73%  '
    'E17' = '  -4.27%  '
    'D18' = '43.968.52'
    'E18' = '  -0.05%  '
    'D19' = '0.0₃0963'
    'E19' = '  -0.57%  '
    'E20' = '  -6.82%  '
    'E21' = '  -1.48%  '
    'E22' = '  +0.47%  '
    'E24' = '  -4.68%  '
    'E25' = '  -4.06%  '
    'E26' = '  +0.25%  '
    'E27' = '  +4.33%  '
    'E28' = '  -4.84%  '
    'E29' = '  -0.62%  '
    'E30' = '  -0.23%  '
    'E31' = '  -0.40%  '
    'E33' = '  -4.76%  '
    'E34' = '  -3.07%  '
    'E35' = '  +0.41%  '
    'E36' = '  -1.12%  '
    'E37' = '  +1.55%  '
    'E38' = '  -7.94%  '
    'E39' = '  -4.81%  '
    'E40' = '  -6.93%  '
    'E41' = '  -6.28%  '
    'E42' = '  -5.03%  '
    'E43' = '  +0.20%  '
    'D44' = '1.747.85'
    'E44' = '  +0.64%  '
    'E45' = '  +3.48%  '
    'E46' = '  -4.39%  '
    'E47' = '  -2.04%  '
    'E48' = '  -3.54%  '
    'B49' = 'EnergySwap'
    'C49' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'E49' = '  +4.44%  '
    'E50' = '  -6.37%  '
    'B51' = 'FraxShare'
    'C51' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'E51' = '  -1.11%  '
}

# Price updates that look like plain numbers (e.g. "306.16", "0.110").
# The source workbook stores these as TEXT (inlineStr), not numbers, so a bare
# assignment would let Excel coerce them to doubles and drop meaningful trailing
# zeros (e.g. "0.110" -> 0.11). A leading apostrophe forces a text literal, and the
# implicit @ (text) style Excel applies for that is reset back to Normal afterwards
# so the cell keeps matching the workbook's other (unstyled) data cells.
$numericTextUpdates = [ordered]@{
    'D5' = '306.16'
    'D6' = '95.19'
    'D9' = '0.524'
    'D10' = '34.84'
    'D12' = '7.22'
    'D16' = '0.824'
    'D17' = '13.52'
    'D20' = '12.22'
    'D21' = '6.31'
    'D22' = '65.54'
    'D23' = '236.54'
    'D24' = '2.92'
    'D25' = '1.97'
    'D27' = '38.56'
    'D28' = '9.88'
    'D29' = '2.20'
    'D30' = '6.09'
    'D31' = '19.96'
    'D32' = '150.53'
    'D33' = '0.0801'
    'D35' = '3.19'
    'D36' = '0.110'
    'D38' = '1.77'
    'D39' = '15.20'
    'D40' = '3.41'
    'D41' = '3.84'
    'D42' = '0.0297'
    'D45' = '84.25'
    'D46' = '0.188'
    'D47' = '100.05'
    'D48' = '4.95'
    'D49' = '14.61'
    'D50' = '69.24'
    'D51' = '8.08'
}

foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

foreach ($ref in $numericTextUpdates.Keys) {
    $ws.Range($ref).Value = "'" + $numericTextUpdates[$ref]
}

foreach ($ref in $numericTextUpdates.Keys) {
    $ws.Range($ref).Style = 'Normal'
}
